$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '29.454.00'
Set-TextValue 'D3' '1.871.82'
Set-TextValue 'E3' '  -0.41%  '
Set-TextValue 'D4' '1.001'
Set-TextValue 'E4' '  -0.16%  '
Set-TextValue 'D5' '243.72'
Set-TextValue 'E5' '  +0.26%  '
Set-TextValue 'D6' '0.7052'
Set-TextValue 'E6' '  -2.21%  '
Set-TextValue 'D7' '1.001'
Set-TextValue 'E7' '  -0.16%  '
Set-TextValue 'D8' '0.07937'
Set-TextValue 'E8' '  -1.02%  '
Set-TextValue 'E9' '  +0.16%  '
Set-TextValue 'E10' '  -1.57%  '
Set-TextValue 'E11' '  -4.19%  '
Set-TextValue 'D12' '1.870.51'
Set-TextValue 'E12' '  -0.79%  '
Set-TextValue 'D13' '93.80'
Set-TextValue 'E13' '  -0.81%  '
Set-TextValue 'D14' '5.173'
Set-TextValue 'E14' '  -1.07%  '
Set-TextValue 'D15' '0.7030'
Set-TextValue 'E15' '  -1.18%  '
Set-TextValue 'D16' '6.503'
Set-TextValue 'E16' '  +1.20%  '
Set-TextValue 'D17' '0.000008559'
Set-TextValue 'E17' '  +1.04%  '
Set-TextValue 'D18' '29.475.74'
Set-TextValue 'E18' '  +0.34%  '
Set-TextValue 'D19' '253.63'
Set-TextValue 'E19' '  +3.89%  '
Set-TextValue 'D20' '2.137.14'
Set-TextValue 'E20' '  +0.36%  '
Set-TextValue 'D21' '13.10'
Set-TextValue 'E21' '  -1.47%  '
Set-TextValue 'D22' '1.000'
Set-TextValue 'E22' '  -0.10%  '
Set-TextValue 'D23' '7.617'
Set-TextValue 'E23' '  -1.62%  '
Set-TextValue 'D24' '1.001'
Set-TextValue 'E24' '  -0.14%  '
Set-TextValue 'D25' '0.1539'
Set-TextValue 'E25' '  -4.25%  '
Set-TextValue 'D26' '9.017'
Set-TextValue 'E26' '  -0.23%  '
Set-TextValue 'D27' '161.33'
Set-TextValue 'E27' '  -0.85%  '
Set-TextValue 'D28' '18.79'
Set-TextValue 'E28' '  +1.49%  '
Set-TextValue 'D29' '1.544'
Set-TextValue 'E29' '  +2.54%  '
Set-TextValue 'D30' '4.310'
Set-TextValue 'E30' '  -2.05%  '
Set-TextValue 'D31' '4.266'
Set-TextValue 'E31' '  -0.40%  '
Set-TextValue 'D32' '1.204'
Set-TextValue 'E32' '  -2.25%  '
Set-TextValue 'D34' '1.897'
Set-TextValue 'E34' '  -1.98%  '
Set-TextValue 'D35' '0.7605'
Set-TextValue 'E35' '  -0.45%  '
Set-TextValue 'D36' '1.182'
Set-TextValue 'E36' '  +0.42%  '
Set-TextValue 'D37' '2.704'
Set-TextValue 'E37' '  +0.24%  '
$ws.Range('B38').Value = 'Maker'
$ws.Range('C38').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue 'D38' '1.281.62'
Set-TextValue 'E38' '  +1.66%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D39' '0.01876'
Set-TextValue 'E39' '  +0.25%  '
Set-TextValue 'D40' '2.757'
Set-TextValue 'D41' '0.9009'
Set-TextValue 'E41' '  -0.41%  '
Set-TextValue 'D42' '109.91'
Set-TextValue 'E42' '  -2.86%  '
Set-TextValue 'D43' '5.982'
Set-TextValue 'E43' '  -7.08%  '
Set-TextValue 'D44' '70.99'
Set-TextValue 'E44' '  -4.26%  '
Set-TextValue 'D45' '1.000'
Set-TextValue 'E45' '  -0.18%  '
Set-TextValue 'D46' '2.037.34'
Set-TextValue 'E46' '  +0.57%  '
Set-TextValue 'D47' '0.00000000126'
Set-TextValue 'E47' '  -3.52%  '
Set-TextValue 'D48' '9.637'
Set-TextValue 'E48' '  +1.64%  '
Set-TextValue 'E49' '  +0.22%  '
Set-TextValue 'D50' '0.5171'
Set-TextValue 'E50' '  -0.47%  '
Set-TextValue 'D51' '0.4310'
Set-TextValue 'E51' '  -0.65%  '
